$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

$ws.Range("C3:C17").Value = "Y"

$ws.Range("C17").Select()
